# Revert the "1-4" / "5-8" / "3-8" range labels back to the "1..4" / "5..8" / "3..8"
# double-dot format, rename the first sheet, and restore the selections that were
# left on A1 (sheet1) / K8 (sheet2).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 rename: "sheet1 1" -> "sheet1" ---
$ws1.Name = "sheet1"

# --- Shared-string text fixes: "-" -> ".." ---
$ws1.Range("A22").Value = "1..4"
$ws1.Range("A23").Value = "5..8"
$ws1.Range("A36").Value = "3..8"

$ws2.Range("B25").Value = "1..4"
$ws2.Range("B26").Value = "5..8"
$ws2.Range("B39").Value = "3..8"

# --- Restore selections ---
$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("K8").Select() | Out-Null
